$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2060.6
$ws.Range("I62").Value = 3900
$ws.Range("J62").Value = 1856.2222
$ws.Range("K62").Value = 3900
$ws.Range("L62").Value = 1856.2222
$ws.Range("M62").Value = -3276
$ws.Range("N62").Value = -3104.2222
$ws.Range("H65").Value = 2060.6
$ws.Range("I65").Value = 3900
$ws.Range("J65").Value = 1856.2222
$ws.Range("K65").Value = 19500
$ws.Range("L65").Value = 9281.110999999999
$ws.Range("M65").Value = -16380
$ws.Range("N65").Value = -15521.111
$ws.Range("H116").Value = 3067.5264
$ws.Range("I116").Value = 3620.6
$ws.Range("J116").Value = 2870
$ws.Range("K116").Value = 3620.6
$ws.Range("L116").Value = 2870
$ws.Range("M116").Value = -178.5999999999999
$ws.Range("N116").Value = -9754
$ws.Range("H132").Value = 4172.551
$ws.Range("I132").Value = 4515.2764
$ws.Range("J132").Value = 3440.3635
$ws.Range("K132").Value = 13545.8292
$ws.Range("L132").Value = 10321.0905
$ws.Range("M132").Value = -11015.8292
$ws.Range("N132").Value = -15381.0905
$ws.Range("H137").Value = 3600.42
$ws.Range("I137").Value = 1075.3334
$ws.Range("J137").Value = 5931.269
$ws.Range("K137").Value = 3226.0002
$ws.Range("L137").Value = 17793.807
$ws.Range("M137").Value = -676.0001999999999
$ws.Range("N137").Value = -22893.807

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 276012.5
$ws.Range("I2").Value = 1000011
$ws.Range("K2").Value = 1000011
$ws.Range("M2").Value = -999898
$ws.Range("H45").Value = 1585.8334
$ws.Range("I45").Value = 1239.3334
$ws.Range("J45").Value = 1932.3334
$ws.Range("K45").Value = 1239.3334
$ws.Range("L45").Value = 1932.3334
$ws.Range("M45").Value = -862.3334
$ws.Range("N45").Value = -2686.3334
$ws.Range("H116").Value = 276012.5
$ws.Range("I116").Value = 1000011
$ws.Range("K116").Value = 1000011
$ws.Range("M116").Value = -997717
$ws.Range("H122").Value = 1566.3334
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1566.3334
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 4699.0002
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -9599.0002
$ws.Range("H132").Value = 19683.834
$ws.Range("I132").Value = 20327.334
$ws.Range("J132").Value = 19040.334
$ws.Range("K132").Value = 60982.00199999999
$ws.Range("L132").Value = 57121.00199999999
$ws.Range("M132").Value = -58452.00199999999
$ws.Range("N132").Value = -62181.00199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 276012.5
$ws.Range("I3").Value = 1000011
$ws.Range("K3").Value = 1000011
$ws.Range("M3").Value = -999897
$ws.Range("H105").Value = 2179.3076
$ws.Range("I105").Value = 1373.3334
$ws.Range("J105").Value = 2421.1
$ws.Range("K105").Value = 1373.3334
$ws.Range("L105").Value = 2421.1
$ws.Range("M105").Value = 373.6666
$ws.Range("N105").Value = -5915.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 24394618
$ws.Range("I132").Value = 32263230
$ws.Range("J132").Value = 1920.2
$ws.Range("K132").Value = 96789690
$ws.Range("L132").Value = 5760.6
$ws.Range("M132").Value = -96787160
$ws.Range("N132").Value = -10820.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 871.75
$ws.Range("I132").Value = 871.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7845.75
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5315.75
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3449.2058
$ws.Range("I102").Value = 3629.25
$ws.Range("J102").Value = 3017.1
$ws.Range("K102").Value = 3629.25
$ws.Range("L102").Value = 3017.1
$ws.Range("M102").Value = -2007.25
$ws.Range("N102").Value = -6261.1
$ws.Range("H122").Value = 33938.332
$ws.Range("I122").Value = 100007
$ws.Range("J122").Value = 904
$ws.Range("K122").Value = 300021
$ws.Range("L122").Value = 2712
$ws.Range("M122").Value = -297571
$ws.Range("N122").Value = -7612
$ws.Range("H132").Value = 4322.6514
$ws.Range("I132").Value = 6243.04
$ws.Range("J132").Value = 1655.4445
$ws.Range("K132").Value = 18729.12
$ws.Range("L132").Value = 4966.333500000001
$ws.Range("M132").Value = -16199.12
$ws.Range("N132").Value = -10026.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1391.6285
$ws.Range("I61").Value = 1086.1765
$ws.Range("J61").Value = 1680.1111
$ws.Range("K61").Value = 1086.1765
$ws.Range("L61").Value = 1680.1111
$ws.Range("M61").Value = -884.1765
$ws.Range("N61").Value = -2084.1111
$ws.Range("H113").Value = 1391.6285
$ws.Range("I113").Value = 1086.1765
$ws.Range("J113").Value = 1680.1111
$ws.Range("K113").Value = 1086.1765
$ws.Range("L113").Value = 1680.1111
$ws.Range("M113").Value = 1083.8235
$ws.Range("N113").Value = -6020.1111
$ws.Range("H122").Value = 7079.826
$ws.Range("I122").Value = 10037
$ws.Range("J122").Value = 3235.5
$ws.Range("K122").Value = 30111
$ws.Range("L122").Value = 9706.5
$ws.Range("M122").Value = -27661
$ws.Range("N122").Value = -14606.5
$ws.Range("H132").Value = 13162
$ws.Range("I132").Value = 33150.5
$ws.Range("J132").Value = 4278.222
$ws.Range("K132").Value = 99451.5
$ws.Range("L132").Value = 12834.666
$ws.Range("M132").Value = -96921.5
$ws.Range("N132").Value = -17894.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7409236.5
$ws.Range("I122").Value = 11766350
$ws.Range("J122").Value = 2144
$ws.Range("K122").Value = 35299050
$ws.Range("L122").Value = 6432
$ws.Range("M122").Value = -35296600
$ws.Range("N122").Value = -11332
$ws.Range("H132").Value = 2477.4
$ws.Range("I132").Value = 3356.525
$ws.Range("J132").Value = 1070.8
$ws.Range("K132").Value = 10069.575
$ws.Range("L132").Value = 3212.4
$ws.Range("M132").Value = -8272.4
